$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    # Force the literal text into the cell without Excel re-casting
    # numeric-looking strings (e.g. "311.66") into real numbers, and
    # without leaving a permanent "Text" number-format behind.
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.ClearFormats()
}

# --- Row 47 / 48: Quant and RenderToken swap places in the ranking ---
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D47") "1.948"
$ws.Range("E47").Value = "  +6.89%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "106.10"
$ws.Range("E48").Value = "  -0.29%  "

# --- Per-row Price (D) / Volume(1h) (E) refresh ---
Set-TextValue $ws.Range("D2") "27.021.95"
$ws.Range("E2").Value = "  -1.34%  "
Set-TextValue $ws.Range("D3") "1.825.04"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue $ws.Range("D5") "311.66"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  -0.25%  "
Set-TextValue $ws.Range("D7") "0.4362"
$ws.Range("E7").Value = "  +1.47%  "
Set-TextValue $ws.Range("D8") "0.3676"
$ws.Range("E8").Value = "  -0.70%  "
Set-TextValue $ws.Range("D9") "0.07268"
$ws.Range("E9").Value = "  +0.15%  "
Set-TextValue $ws.Range("D10") "0.8453"
$ws.Range("E10").Value = "  -2.48%  "
Set-TextValue $ws.Range("D11") "20.67"
$ws.Range("E11").Value = "  -2.39%  "
Set-TextValue $ws.Range("D12") "1.824.93"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -0.40%  "
Set-TextValue $ws.Range("D14") "0.07074"
$ws.Range("E14").Value = "  -0.36%  "
Set-TextValue $ws.Range("D15") "5.296"
$ws.Range("E15").Value = "  -1.09%  "
Set-TextValue $ws.Range("D16") "89.57"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("E17").Value = "  -0.30%  "
Set-TextValue $ws.Range("D18") "0.000008776"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("E20").Value = "  -2.21%  "
Set-TextValue $ws.Range("D21") "27.158.88"
$ws.Range("E21").Value = "  -0.89%  "
Set-TextValue $ws.Range("D22") "5.152"
$ws.Range("E22").Value = "  -0.34%  "
Set-TextValue $ws.Range("D23") "10.87"
$ws.Range("E23").Value = "  +0.05%  "
Set-TextValue $ws.Range("D24") "2.056.71"
$ws.Range("E24").Value = "  -0.07%  "
Set-TextValue $ws.Range("D25") "1.993"
$ws.Range("E25").Value = "  -0.98%  "
Set-TextValue $ws.Range("D26") "151.63"
$ws.Range("E26").Value = "  -1.11%  "
Set-TextValue $ws.Range("D27") "2.209"
Set-TextValue $ws.Range("D28") "18.28"
$ws.Range("E28").Value = "  -0.95%  "
Set-TextValue $ws.Range("D29") "5.236"
$ws.Range("E29").Value = "  -1.23%  "
Set-TextValue $ws.Range("D30") "116.85"
$ws.Range("E30").Value = "  -0.42%  "
Set-TextValue $ws.Range("D31") "0.08779"
$ws.Range("E31").Value = "  -0.66%  "
Set-TextValue $ws.Range("D32") "1.178"
$ws.Range("E32").Value = "  -2.01%  "
Set-TextValue $ws.Range("D33") "0.7421"
$ws.Range("E33").Value = "  -3.30%  "
Set-TextValue $ws.Range("D34") "2.907"
$ws.Range("E34").Value = "  +1.41%  "
Set-TextValue $ws.Range("D35") "4.430"
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("E37").Value = "  -2.28%  "
$ws.Range("E38").Value = "  -0.68%  "
Set-TextValue $ws.Range("D39") "0.05231"
$ws.Range("E39").Value = "  -0.91%  "
Set-TextValue $ws.Range("D40") "7.240"
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  +1.03%  "
Set-TextValue $ws.Range("D43") "0.5157"
$ws.Range("E43").Value = "  +1.71%  "
Set-TextValue $ws.Range("D44") "8.565"
$ws.Range("E44").Value = "  -1.16%  "
Set-TextValue $ws.Range("D45") "10.59"
$ws.Range("E45").Value = "  +0.02%  "
Set-TextValue $ws.Range("D46") "0.4787"
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("E49").Value = "  -0.25%  "
Set-TextValue $ws.Range("D50") "0.06347"
$ws.Range("E50").Value = "  -1.18%  "
Set-TextValue $ws.Range("D51") "1.659"
$ws.Range("E51").Value = "  -0.70%  "
